$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "59.404.40"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "2.643.63"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").Formula = "'517.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").Formula = "'145.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Formula = "'0.574"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("D9").Value = "2.652.00"
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").Formula = "'6.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.79%  "
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").Value = "3.106.35"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").Value = "59.382.17"
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").Formula = "'20.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.08%  "
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").Value = "2.649.21"
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("D19").Formula = "'349.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("E20").Value = "  -1.97%  "
$ws.Range("D21").Formula = "'10.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.52%  "
$ws.Range("D22").Formula = "'6.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.26%  "
$ws.Range("D23").Formula = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("D24").Formula = "'62.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.00%  "
$ws.Range("D25").Formula = "'0.416"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.34%  "
$ws.Range("D26").Formula = "'0.166"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.28%  "
$ws.Range("D27").Formula = "'0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("D29").Formula = "'7.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("D30").Formula = "'0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").Formula = "'6.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").Formula = "'18.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").Formula = "'149.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("B35").Value = "SuiNetwork"
$ws.Range("C35").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D35").Formula = "'0.951"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -10.90%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Formula = "'4.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("D37").Formula = "'1.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("D38").Formula = "'0.861"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("D39").Formula = "'36.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("E40").Value = "  +4.28%  "
$ws.Range("D41").Formula = "'3.67"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.10%  "
$ws.Range("E42").Value = "  +0.36%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Formula = "'277.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.72%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Formula = "'0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").Formula = "'0.601"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.92%  "
$ws.Range("D46").Formula = "'19.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("D47").Value = "2.084.24"
$ws.Range("E47").Value = "  +4.77%  "
$ws.Range("D48").Formula = "'0.0529"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.60%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").Formula = "'10.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.55%  "
$ws.Range("D51").Formula = "'4.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.86%  "
